# Insert a new data row at row 110 (pushing existing rows 110..231 down to
# 111..232), then populate the freshly inserted row with the new record.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Rows.Item(110).Insert()

$ws.Cells.Item(110, 1).Value = 10
$ws.Cells.Item(110, 2).Value = "Vega Modelo de Temuco"
$ws.Cells.Item(110, 3).Value = "La Araucanía"
$ws.Cells.Item(110, 4).Value = 44539
$ws.Cells.Item(110, 5).Value = 9
$ws.Cells.Item(110, 6).Value = 100112009
$ws.Cells.Item(110, 7).Value = "Acelga"
$ws.Cells.Item(110, 8).Value = "Sin especificar"
$ws.Cells.Item(110, 9).Value = "Primera"
$ws.Cells.Item(110, 10).Value = 115
$ws.Cells.Item(110, 11).Value = 8000
$ws.Cells.Item(110, 12).Value = 9000
$ws.Cells.Item(110, 13).Value = 8565
$ws.Cells.Item(110, 14).Value = "`$/docena de atados (12 kilos)"
$ws.Cells.Item(110, 15).Value = "Provincia de Cautín"
$ws.Cells.Item(110, 16).Value = 714
$ws.Cells.Item(110, 17).Value = 12
$ws.Cells.Item(110, 18).Value = "Hortaliza"
